$wb = $excel.ActiveWorkbook

# --- Sheet 1: Generator Data ---
$ws1 = $wb.Worksheets.Item("Generator Data")
$ws1.Range("B2").Value = 455432.19341
$ws1.Range("B3").Value = 191327.064451541
$ws1.Range("B4").Value = 19132.7064451541
$ws1.Range("B5").Value = 3170223.58785

# --- Sheet 2: Yearly Fuel Costs ---
$ws2 = $wb.Worksheets.Item("Yearly Fuel Costs")

# Delete rows 7 through 21 (only rows 1-6 remain)
$ws2.Range("A7:B21").EntireRow.Delete()

# Set the new values for the remaining data rows
$ws2.Range("B2").Value = 706126.2457347935
$ws2.Range("B3").Value = 706126.9053339667
$ws2.Range("B4").Value = 706127.4380772213
$ws2.Range("B5").Value = 706128.0094528636
$ws2.Range("B6").Value = 706128.4838875127
